# transitioned to shout framework and implemented out-of-ammo shout
#
# This moves the data that used to live in column B (rows 13-15) over to
# column C (rows 13-15), applying the same "highlighted" style (s=7 /
# fillId=2, the style originally only used by B13) to all three of them.
# Since nothing is left using the green solid fill afterwards, the fill,
# its conditional-formatting dxf, and the conditional formatting rule
# that depended on column B being filled in are all cleaned up too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values currently sitting in B13:B15 before we start moving
# things around.
$b13 = $ws.Range("B13").Value
$b14 = $ws.Range("B14").Value
$b15 = $ws.Range("B15").Value

# B13 already uses the "highlighted" style; re-use it for the cells we are
# about to populate in column C, then clear out column B.
$ws.Range("C13").Value = $b13
$ws.Range("C14").Value = $b14
$ws.Range("C15").Value = $b15

$ws.Range("C13:C15").Style = $ws.Range("B13").Style

$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()

# The old conditional formatting rule (highlighting non-blank cells in
# B2:B39 using the green fill) is no longer needed now that column B is
# empty of data again, so remove it along with the associated format.
$ws.Range("B2:B39").FormatConditions.Delete()

# Update the active selection to match where the edit was made.
$ws.Range("C13").Select()
